$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$rng = $ws3.Range("A1:F1")
$rng.Font.Bold = $true
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160
